$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A22").Value = "22TRD1111"
$ws.Range("B22").Value = "Note"
$ws.Range("C22").Value = "KUDELA"
$ws.Range("D22").Value = "JUSTIN"
$ws.Range("E22").Value = "FIRE"
$ws.Range("F22").Value = 2244.22
$ws.Range("G22").Value = "MM"
$ws.Range("H22").Value = "Y"

$ws.Range("A23").Value = "24TRD2222"
$ws.Range("B23").Value = "Note"
$ws.Range("C23").Value = "SMITH"
$ws.Range("D23").Value = "MAXIMUS"
$ws.Range("E23").Value = "WATER"
$ws.Range("F23").Value = 4422.23
$ws.Range("G23").Value = "UCM"
$ws.Range("H23").Value = "N"

$ws.Range("I23").Select()
